$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2543.1667
$ws.Range("I40").Value = 2641.818
$ws.Range("J40").Value = 2459.6924
$ws.Range("K40").Value = 2641.818
$ws.Range("L40").Value = 2459.6924
$ws.Range("M40").Value = -2466.818
$ws.Range("N40").Value = -2809.6924
# Row 76
$ws.Range("H76").Value = 3624.054
$ws.Range("I76").Value = 3002.8518
$ws.Range("J76").Value = 5301.3
$ws.Range("K76").Value = 3002.8518
$ws.Range("L76").Value = 5301.3
$ws.Range("M76").Value = -2687.8518
$ws.Range("N76").Value = -5931.3
# Row 79
$ws.Range("H79").Value = 3624.054
$ws.Range("I79").Value = 3002.8518
$ws.Range("J79").Value = 5301.3
$ws.Range("K79").Value = 3002.8518
$ws.Range("L79").Value = 5301.3
$ws.Range("M79").Value = -1910.8518
$ws.Range("N79").Value = -7485.3
# Row 112
$ws.Range("H112").Value = 1225.3914
$ws.Range("J112").Value = 1256.381
$ws.Range("L112").Value = 3769.143
$ws.Range("N112").Value = -5985.143
# Row 129
$ws.Range("H129").Value = 1007.96155
$ws.Range("I129").Value = 512.125
$ws.Range("J129").Value = 1098.1136
$ws.Range("K129").Value = 1536.375
$ws.Range("L129").Value = 3294.3408
$ws.Range("M129").Value = 3463.625
$ws.Range("N129").Value = -13294.3408
# Row 132
$ws.Range("H132").Value = 27676.334
$ws.Range("I132").Value = 38856.79
$ws.Range("J132").Value = 3079.3333
$ws.Range("K132").Value = 116570.37
$ws.Range("L132").Value = 9237.999899999999
$ws.Range("M132").Value = -114040.37
$ws.Range("N132").Value = -14297.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3351.91
$ws.Range("I32").Value = 3351.91
$ws.Range("K32").Value = 3351.91
$ws.Range("M32").Value = -3064.91
# Row 35
$ws.Range("H35").Value = 2350
$ws.Range("I35").Value = 2350
$ws.Range("K35").Value = 2350
$ws.Range("M35").Value = -1944

$ws = $wb.Worksheets.Item("BSM")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 99
$ws.Range("H99").Value = 1515.2142
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1554.8462
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1554.8462
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4550.8462
# Row 105
$ws.Range("H105").Value = 3057.9
$ws.Range("I105").Value = 2951.2856
$ws.Range("J105").Value = 3306.6667
$ws.Range("K105").Value = 2951.2856
$ws.Range("L105").Value = 3306.6667
$ws.Range("M105").Value = -1204.2856
$ws.Range("N105").Value = -6800.6667
# Row 107
$ws.Range("H107").Value = 687.3333
$ws.Range("I107").Value = 679.2308
$ws.Range("J107").Value = 740
$ws.Range("K107").Value = 679.2308
$ws.Range("L107").Value = 740
$ws.Range("M107").Value = 1240.7692
$ws.Range("N107").Value = -4580

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 29623.62
$ws.Range("I51").Value = 7500
$ws.Range("J51").Value = 31952.422
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 31952.422
$ws.Range("M51").Value = -6764
$ws.Range("N51").Value = -33424.422
# Row 61
$ws.Range("H61").Value = 29623.62
$ws.Range("I61").Value = 7500
$ws.Range("J61").Value = 31952.422
$ws.Range("K61").Value = 7500
$ws.Range("L61").Value = 31952.422
$ws.Range("M61").Value = -7152
$ws.Range("N61").Value = -32648.422
# Row 99
$ws.Range("H99").Value = 17243.143
$ws.Range("I99").Value = 2300.5
$ws.Range("J99").Value = 37166.668
$ws.Range("K99").Value = 2300.5
$ws.Range("L99").Value = 37166.668
$ws.Range("M99").Value = -802.5
$ws.Range("N99").Value = -40162.668
# Row 126
$ws.Range("H126").Value = 17243.143
$ws.Range("I126").Value = 2300.5
$ws.Range("J126").Value = 37166.668
$ws.Range("K126").Value = 6901.5
$ws.Range("L126").Value = 111500.004
$ws.Range("M126").Value = -4431.5
$ws.Range("N126").Value = -116440.004
# Row 132
$ws.Range("H132").Value = 596855.0600000001
$ws.Range("I132").Value = 676808.4399999999
$ws.Range("J132").Value = 5199.8
$ws.Range("K132").Value = 2030425.32
$ws.Range("L132").Value = 15599.4
$ws.Range("M132").Value = -2027895.32
$ws.Range("N132").Value = -20659.4
# Row 141
$ws.Range("H141").Value = 38761.668
$ws.Range("I141").Value = 20100
$ws.Range("J141").Value = 43427.082
$ws.Range("K141").Value = 20100
$ws.Range("L141").Value = 43427.082
$ws.Range("M141").Value = -14920
$ws.Range("N141").Value = -53787.082

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 950.7
$ws.Range("I98").Value = 1029
$ws.Range("J98").Value = 768
$ws.Range("K98").Value = 3087
$ws.Range("L98").Value = 2304
$ws.Range("M98").Value = -1589
$ws.Range("N98").Value = -5300
# Row 131
$ws.Range("H131").Value = 858.6799999999999
$ws.Range("I131").Value = 275
$ws.Range("J131").Value = 883
$ws.Range("K131").Value = 825
$ws.Range("L131").Value = 2649
$ws.Range("M131").Value = 4215
$ws.Range("N131").Value = -12729

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 32500
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 36000
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 36000
$ws.Range("M24").Value = -827
$ws.Range("N24").Value = -36346
# Row 70
$ws.Range("H70").Value = 5659.591
$ws.Range("I70").Value = 4955.5835
$ws.Range("J70").Value = 6504.4
$ws.Range("K70").Value = 4955.5835
$ws.Range("L70").Value = 6504.4
$ws.Range("M70").Value = -4685.5835
$ws.Range("N70").Value = -7044.4
# Row 73
$ws.Range("H73").Value = 5659.591
$ws.Range("I73").Value = 4955.5835
$ws.Range("J73").Value = 6504.4
$ws.Range("K73").Value = 4955.5835
$ws.Range("L73").Value = 6504.4
$ws.Range("M73").Value = -4019.5835
$ws.Range("N73").Value = -8376.4
# Row 80
$ws.Range("H80").Value = 3161.2
$ws.Range("I80").Value = 2933.3333
$ws.Range("J80").Value = 3503
$ws.Range("K80").Value = 2933.3333
$ws.Range("L80").Value = 3503
$ws.Range("M80").Value = -1935.3333
$ws.Range("N80").Value = -5499
# Row 83
$ws.Range("H83").Value = 3161.2
$ws.Range("I83").Value = 2933.3333
$ws.Range("J83").Value = 3503
$ws.Range("K83").Value = 14666.6665
$ws.Range("L83").Value = 17515
$ws.Range("M83").Value = -9674.666499999999
$ws.Range("N83").Value = -27499
# Row 122
$ws.Range("H122").Value = 72372.64
$ws.Range("I122").Value = 95162.03999999999
$ws.Range("J122").Value = 4004.4443
$ws.Range("K122").Value = 285486.12
$ws.Range("L122").Value = 12013.3329
$ws.Range("M122").Value = -283036.12
$ws.Range("N122").Value = -16913.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 922.0454999999999
$ws.Range("I46").Value = 747.73334
$ws.Range("J46").Value = 1295.5714
$ws.Range("K46").Value = 747.73334
$ws.Range("L46").Value = 1295.5714
$ws.Range("M46").Value = -559.73334
$ws.Range("N46").Value = -1671.5714
# Row 122
$ws.Range("H122").Value = 6538224.5
$ws.Range("I122").Value = 15874791
$ws.Range("J122").Value = 2628
$ws.Range("K122").Value = 47624373
$ws.Range("L122").Value = 7884
$ws.Range("M122").Value = -47621923
$ws.Range("N122").Value = -12784

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 331.88
$ws.Range("I107").Value = 263.1875
$ws.Range("J107").Value = 454
$ws.Range("K107").Value = 789.5625
$ws.Range("L107").Value = 1362
$ws.Range("M107").Value = 1130.4375
$ws.Range("N107").Value = -5202
